$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "...payment processing to|  handl|e| charitable donations" (4 runs) ->
# "...payment processing to handle charitable donations" (merged run)
$found1 = $d.Content.Find.Execute(
    "to handle charitable donations", $true, $false, $false, $false, $false,
    $true, 1, $false, "to handle charitable donations", 2)

# --- Change 2 -------------------------------------------------------------
# "Angular 6 UI enhancement..." -> "Angular 4 UI enhancement..."
$found2 = $d.Content.Find.Execute(
    "Angular 6 UI enhancement", $true, $false, $false, $false, $false,
    $true, 1, $false, "Angular 4 UI enhancement", 2)

# --- Change 3 -------------------------------------------------------------
# "...AngularJS, Node|JS|, Bootstrap, SharePoint API..." (3 runs) ->
# "...AngularJS, NodeJS, Bootstrap, SharePoint API..." (merged run)
$found3 = $d.Content.Find.Execute(
    "AngularJS, NodeJS, Bootstrap, SharePoint API", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "AngularJS, NodeJS, Bootstrap, SharePoint API", 2)

Write-Host "Change1:" $found1 "Change2:" $found2 "Change3:" $found3
